# Auto-generated edit script: update cryptos list (prices & 1h volume %)
# to the refreshed GitHub Actions snapshot, plus the Cronos/Quant row swap.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '23.108.70'
$ws.Range("E2").Value = '  -3.42%  '
$ws.Range("D3").Value = '1.597.60'
$ws.Range("E3").Value = '  -3.12%  '
$ws.Range("D4").Value = '1.003'
$ws.Range("E4").Value = '  +0.35%  '
$ws.Range("D5").Value = '1.002'
$ws.Range("E5").Value = '  +0.19%  '
$ws.Range("D6").Value = '301.72'
$ws.Range("E6").Value = '  -2.43%  '
$ws.Range("E7").Value = '  -3.26%  '
$ws.Range("D8").Value = '0.3656'
$ws.Range("E8").Value = '  -4.44%  '
$ws.Range("D9").Value = '47.82'
$ws.Range("E9").Value = '  -6.78%  '
$ws.Range("D10").Value = '1.003'
$ws.Range("E10").Value = '  +0.32%  '
$ws.Range("D11").Value = '1.276'
$ws.Range("E11").Value = '  -5.34%  '
$ws.Range("D12").Value = '0.08079'
$ws.Range("E12").Value = '  -4.17%  '
$ws.Range("D13").Value = '22.95'
$ws.Range("E13").Value = '  -3.77%  '
$ws.Range("D14").Value = '6.627'
$ws.Range("E14").Value = '  -6.70%  '
$ws.Range("D15").Value = '7.619'
$ws.Range("E15").Value = '  -2.40%  '
$ws.Range("D16").Value = '0.00001266'
$ws.Range("E16").Value = '  -3.39%  '
$ws.Range("D17").Value = '1.596.27'
$ws.Range("E17").Value = '  -3.34%  '
$ws.Range("D19").Value = '0.06788'
$ws.Range("E19").Value = '  -2.78%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '18.40'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  -6.64%  '
$ws.Range("D21").Value = '6.583'
$ws.Range("E21").Value = '  -4.08%  '
$ws.Range("E22").Value = '  +0.17%  '
$ws.Range("D23").Value = '13.04'
$ws.Range("E23").Value = '  -3.87%  '
$ws.Range("D24").Value = '23.122.42'
$ws.Range("E24").Value = '  -3.30%  '
$ws.Range("D25").Value = '2.366'
$ws.Range("E25").Value = '  -4.88%  '
$ws.Range("D26").Value = '2.894'
$ws.Range("E26").Value = '  -2.94%  '
$ws.Range("D27").Value = '21.08'
$ws.Range("E27").Value = '  -4.20%  '
$ws.Range("D28").Value = '150.98'
$ws.Range("E28").Value = '  -0.82%  '
$ws.Range("D29").Value = '5.244'
$ws.Range("E29").Value = '  -3.42%  '
$ws.Range("D30").Value = '131.98'
$ws.Range("E30").Value = '  -5.07%  '
$ws.Range("D31").Value = '2.448'
$ws.Range("E31").Value = '  -1.50%  '
$ws.Range("D32").Value = '7.069'
$ws.Range("E32").Value = '  -8.65%  '
$ws.Range("D33").Value = '1.772.50'
$ws.Range("E33").Value = '  -3.21%  '
$ws.Range("D34").Value = '0.9807'
$ws.Range("E34").Value = '  -4.41%  '
$ws.Range("E35").Value = '  -3.78%  '
$ws.Range("D36").Value = '0.02779'
$ws.Range("E36").Value = '  -5.94%  '
$ws.Range("D37").Value = '6.304'
$ws.Range("E37").Value = '  -5.81%  '
$ws.Range("D38").Value = '0.2542'
$ws.Range("E38").Value = '  -5.19%  '
$ws.Range("D39").Value = '0.08866'
$ws.Range("E39").Value = '  -2.56%  '
$ws.Range("D40").Value = '10.05'
$ws.Range("E40").Value = '  -6.69%  '
$ws.Range("D41").Value = '1.397'
$ws.Range("E41").Value = '  -2.36%  '
$ws.Range("D42").Value = '0.7152'
$ws.Range("E42").Value = '  -5.07%  '
$ws.Range("D43").Value = '12.76'
$ws.Range("E43").Value = '  -5.08%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '16.00'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  -0.61%  '
$ws.Range("D45").Value = '0.6634'
$ws.Range("E45").Value = '  -3.97%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '2.310'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  -5.53%  '
$ws.Range("D47").Value = '1.001'
$ws.Range("E47").Value = '  +0.08%  '
$ws.Range("E48").Value = '  -2.58%  '
$ws.Range("B49").Value = 'Quant'
$ws.Range("C49").Value = 'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt'
$ws.Range("D49").Value = '132.03'
$ws.Range("E49").Value = '  -1.71%  '
$ws.Range("B50").Value = 'Cronos'
$ws.Range("C50").Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range("D50").Value = '0.07976'
$ws.Range("E50").Value = '  -3.95%  '
$ws.Range("D51").Value = '1.172'
$ws.Range("E51").Value = '  -3.94%  '
